$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados")
$ws.Columns.Item(14).Delete()
